# Case_0_19 res_line/pl_mw.xlsx update ("case with 380 kV done")
# Rewrites the line-loading result table (rows 2:25) for columns C:I and K:L
# with the recomputed values from the 380 kV case run. Column J (all zeros)
# and columns A, B, M, N, O are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Columns C:I, rows 2:25 ---------------------------------------------
$blockCI = New-Object 'object[,]' 24,7
$blockCI[0,0] = 0.187975658643623
$blockCI[0,1] = 0.04894305529089493
$blockCI[0,2] = 0.1344135178954993
$blockCI[0,3] = 1.141000391209374
$blockCI[0,4] = 1.006847600612446
$blockCI[0,5] = 1.019962940504101
$blockCI[0,6] = 1.131262883943243
$blockCI[1,0] = 0.1823087664189842
$blockCI[1,1] = 0.04868370703032809
$blockCI[1,2] = 0.1308466779653514
$blockCI[1,3] = 1.145454231633074
$blockCI[1,4] = 1.013717250105017
$blockCI[1,5] = 1.030586693653234
$blockCI[1,6] = 1.136932692954822
$blockCI[2,0] = 0.1788923657700394
$blockCI[2,1] = 0.04853878181821258
$blockCI[2,2] = 0.1287250105586075
$blockCI[2,3] = 1.149235787768241
$blockCI[2,4] = 1.018985421664098
$blockCI[2,5] = 1.037848361948718
$blockCI[2,6] = 1.141375619945109
$blockCI[3,0] = 0.1775160746595645
$blockCI[3,1] = 0.04848333679548844
$blockCI[3,2] = 0.1278775951902986
$blockCI[3,3] = 1.15103905357784
$blockCI[3,4] = 1.021395089903862
$blockCI[3,5] = 1.040992859409414
$blockCI[3,6] = 1.143427150055693
$blockCI[4,0] = 0.1772885057303597
$blockCI[4,1] = 0.0484743488291528
$blockCI[4,2] = 0.1277379198098636
$blockCI[4,3] = 1.151354294680644
$blockCI[4,4] = 1.021811052689586
$blockCI[4,5] = 1.041526180475898
$blockCI[4,6] = 1.143782337231656
$blockCI[5,0] = 0.1788737400921718
$blockCI[5,1] = 0.04853801941941427
$blockCI[5,2] = 0.1287135124704939
$blockCI[5,3] = 1.14925904672485
$blockCI[5,4] = 1.019016856679599
$blockCI[5,5] = 1.037890020010153
$blockCI[5,6] = 1.141402312875165
$blockCI[6,0] = 0.1860086490434014
$blockCI[6,1] = 0.04885066839242569
$blockCI[6,2] = 0.1331694670379662
$blockCI[6,3] = 1.142318193524069
$blockCI[6,4] = 1.008997548936136
$blockCI[6,5] = 1.023472447383924
$blockCI[6,6] = 1.133017754794409
$blockCI[7,0] = 0.20049936208774
$blockCI[7,1] = 0.04957688409280081
$blockCI[7,2] = 0.1424515359767256
$blockCI[7,3] = 1.137059962545237
$blockCI[7,4] = 0.9977387732792806
$blockCI[7,5] = 1.001079995551805
$blockCI[7,6] = 1.124243583198449
$blockCI[8,0] = 0.2114494559152149
$blockCI[8,1] = 0.05017891554908971
$blockCI[8,2] = 0.1496054836867771
$blockCI[8,3] = 1.138354748893875
$blockCI[8,4] = 0.9946604034117854
$blockCI[8,5] = 0.9882405648529016
$blockCI[8,6] = 1.122525652085855
$blockCI[9,0] = 0.2164969004055308
$blockCI[9,1] = 0.05046759049043459
$blockCI[9,2] = 0.1529332925827731
$blockCI[9,3] = 1.140078150508117
$blockCI[9,4] = 0.9944048614466681
$blockCI[9,5] = 0.9831899402261541
$blockCI[9,6] = 1.122782612589134
$blockCI[10,0] = 0.2184177289688023
$blockCI[10,1] = 0.05057902660885105
$blockCI[10,2] = 0.1542040457126959
$blockCI[10,3] = 1.140895019662054
$blockCI[10,4] = 0.9944740987722014
$blockCI[10,5] = 0.9813915082108338
$blockCI[10,6] = 1.123030179069204
$blockCI[11,0] = 0.2180036236845524
$blockCI[11,1] = 0.05055493262563715
$blockCI[11,2] = 0.1539298951760415
$blockCI[11,3] = 1.140711768991622
$blockCI[11,4] = 0.9944517820972152
$blockCI[11,5] = 0.9817737485336693
$blockCI[11,6] = 1.122970163405455
$blockCI[12,0] = 0.2166547385820081
$blockCI[12,1] = 0.05047671593841585
$blockCI[12,2] = 0.1530376259487696
$blockCI[12,3] = 1.14014205644115
$blockCI[12,4] = 0.9944072227199285
$blockCI[12,5] = 0.9830396918496831
$blockCI[12,6] = 1.122799963222249
$blockCI[13,0] = 0.2158297395207569
$blockCI[13,1] = 0.05042908197333418
$blockCI[13,2] = 0.1524924643576711
$blockCI[13,3] = 1.139814515401582
$blockCI[13,4] = 0.9944015886927531
$blockCI[13,5] = 0.9838299970971178
$blockCI[13,6] = 1.12271530686354
$blockCI[14,0] = 0.2111209213991856
$blockCI[14,1] = 0.05016034721402463
$blockCI[14,2] = 0.1493894833534668
$blockCI[14,3] = 1.138265043981633
$blockCI[14,4] = 0.994700258449285
$blockCI[14,5] = 0.988586577776033
$blockCI[14,6] = 1.122529835122528
$blockCI[15,0] = 0.2082491329724689
$blockCI[15,1] = 0.04999927366703361
$blockCI[15,2] = 0.1475047314892279
$blockCI[15,3] = 1.137605842237548
$blockCI[15,4] = 0.9951776811077053
$blockCI[15,5] = 0.9917073142396617
$blockCI[15,6] = 1.122682689960818
$blockCI[16,0] = 0.2066035892350158
$blockCI[16,1] = 0.04990802263182559
$blockCI[16,2] = 0.1464275833363402
$blockCI[16,3] = 1.137333393933275
$blockCI[16,4] = 0.9955599554776029
$blockCI[16,5] = 0.9935766160399737
$blockCI[16,6] = 1.122868296483318
$blockCI[17,0] = 0.2060475084024063
$blockCI[17,1] = 0.04987736632060091
$blockCI[17,2] = 0.1460640658077494
$blockCI[17,3] = 1.137259440294969
$blockCI[17,4] = 0.9957078379985518
$blockCI[17,5] = 0.994222281257251
$blockCI[17,6] = 1.122947889387909
$blockCI[18,0] = 0.2085541950089009
$blockCI[18,1] = 0.05001627598906566
$blockCI[18,2] = 0.1477046510869471
$blockCI[18,3] = 1.137664963465355
$blockCI[18,4] = 0.9951157048876809
$blockCI[18,5] = 0.9913674092504436
$blockCI[18,6] = 1.122656301200621
$blockCI[19,0] = 0.217050682165052
$blockCI[19,1] = 0.050499632567508
$blockCI[19,2] = 0.1532994196089632
$blockCI[19,3] = 1.140304927991423
$blockCI[19,4] = 0.9944157945195968
$blockCI[19,5] = 0.982664752039426
$blockCI[19,6] = 1.122845869648316
$blockCI[20,0] = 0.2226588105664291
$blockCI[20,1] = 0.0508278948010954
$blockCI[20,2] = 0.1570176282272868
$blockCI[20,3] = 1.142988224326942
$blockCI[20,4] = 0.9949266137114563
$blockCI[20,5] = 0.9776425399421385
$blockCI[20,6] = 1.123846035795268
$blockCI[21,0] = 0.2196606128287044
$blockCI[21,1] = 0.05065156650765346
$blockCI[21,2] = 0.1550274964213614
$blockCI[21,3] = 1.141468074078929
$blockCI[21,4] = 0.9945649274622781
$blockCI[21,5] = 0.9802619269705986
$blockCI[21,6] = 1.123231740280474
$blockCI[22,0] = 0.2084162593723988
$blockCI[22,1] = 0.05000858502614847
$blockCI[22,2] = 0.147614247439698
$blockCI[22,3] = 1.137637902995834
$blockCI[22,4] = 0.9951433886623704
$blockCI[22,5] = 0.9915208462507223
$blockCI[22,6] = 1.122667927215552
$blockCI[23,0] = 0.1965258633609182
$blockCI[23,1] = 0.04936836683871348
$blockCI[23,2] = 0.1398819706826941
$blockCI[23,3] = 1.137581426945033
$blockCI[23,4] = 0.9998782826928476
$blockCI[23,5] = 1.006505382199904
$blockCI[23,6] = 1.125790767114005
$ws.Range("C2:I25").Value = $blockCI

# --- Columns K:L, rows 2:25 ---------------------------------------------
$blockKL = New-Object 'object[,]' 24,2
$blockKL[0,0] = 2.051353989771656
$blockKL[0,1] = 0.1763426918895163
$blockKL[1,0] = 1.843690623782379
$blockKL[1,1] = 0.1720506178564349
$blockKL[2,0] = 1.716211987121824
$blockKL[2,1] = 0.1695250733299645
$blockKL[3,0] = 1.664272126068113
$blockKL[3,1] = 0.1685234280835388
$blockKL[4,0] = 1.655648132605052
$blockKL[4,1] = 0.168358766454638
$blockKL[5,0] = 1.715511470256217
$blockKL[5,1] = 0.1695114534054056
$blockKL[6,0] = 1.979746949826676
$blockKL[6,1] = 0.1748399460303602
$blockKL[7,0] = 2.498079835246699
$blockKL[7,1] = 0.1861646523769309
$blockKL[8,0] = 2.878976170875148
$blockKL[8,1] = 0.1950257748280251
$blockKL[9,0] = 3.052270830510452
$blockKL[9,1] = 0.1991759921722434
$blockKL[10,0] = 3.117895372748649
$blockKL[10,1] = 0.2007648289836652
$blockKL[11,0] = 3.103761907818694
$blockKL[11,1] = 0.200421876336776
$blockKL[12,0] = 3.057669778183993
$blockKL[12,1] = 0.1993063606947771
$blockKL[13,0] = 3.029437179159856
$blockKL[13,1] = 0.1986253230057855
$blockKL[14,0] = 2.86765127410672
$blockKL[14,1] = 0.1947569557610223
$blockKL[15,0] = 2.768405529410245
$blockKL[15,1] = 0.1924144449815373
$blockKL[16,0] = 2.711324264562052
$blockKL[16,1] = 0.1910783133780285
$blockKL[17,0] = 2.691997970304612
$blockKL[17,1] = 0.1906278462124362
$blockKL[18,0] = 2.778970186634695
$blockKL[18,1] = 0.192662647505017
$blockKL[19,0] = 3.071208123272299
$blockKL[19,1] = 0.199633546063481
$blockKL[20,0] = 3.262211185474825
$blockKL[20,1] = 0.2042899603186754
$blockKL[21,0] = 3.160269005228542
$blockKL[21,1] = 0.2017955152659425
$blockKL[22,0] = 2.774193978381334
$blockKL[22,1] = 0.1925504021069031
$blockKL[23,0] = 2.357841927717573
$blockKL[23,1] = 0.1830065565206098
$ws.Range("K2:L25").Value = $blockKL

